$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COS10001-Students")

$ws.Range("H1").Value = "campus"
$ws.Range("H2").Value = "B"

$ws.Range("H2").Select()
